$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsSchedule.Range("E4").Value = 521.2471875
$wsSchedule.Range("F4").Value = 30.64357363315697
$wsSchedule.Range("E5").Value = 517.9304422499999
$wsSchedule.Range("F5").Value = 15.224292835097
$wsDetailed.Range("B35").Value = 59.97758
$wsDetailed.Range("B36").Value = 59.47856
$wsDetailed.Range("B37").Value = 29.02535
$wsDetailed.Range("C37").Value = "historical"
$wsDetailed.Range("B38").Value = 36.48286
$wsDetailed.Range("C38").Value = "historical"
$wsDetailed.Range("B39").Value = 65.56614
$wsDetailed.Range("B40").Value = 80.02
$wsDetailed.Range("B41").Value = 80.08
$wsDetailed.Range("B42").Value = 84.84780000000001
$wsDetailed.Range("B43").Value = 108.26
$wsDetailed.Range("B44").Value = 90.45925
$wsDetailed.Range("B45").Value = 101.25
$wsDetailed.Range("B47").Value = 69.09292000000001
$wsDetailed.Range("B48").Value = 57.47994
$wsDetailed.Range("B49").Value = 59.46624
$wsDetailed.Range("B50").Value = 59.31501
$wsDetailed.Range("B51").Value = 58.20334
$wsDetailed.Range("B52").Value = 57.89078
$wsDetailed.Range("B53").Value = 59.20427
$wsDetailed.Range("B54").Value = 56.98
$wsDetailed.Range("B56").Value = 49.83939
$wsDetailed.Range("B57").Value = 50.15336
$wsDetailed.Range("B58").Value = 56.98
$wsDetailed.Range("B59").Value = 64.97186000000001
$wsDetailed.Range("B60").Value = 65
$wsDetailed.Range("B61").Value = 75.85464
$wsDetailed.Range("B62").Value = 69.90389999999999
$wsDetailed.Range("B65").Value = 36.06
$wsDetailed.Range("B66").Value = 36.07
$wsDetailed.Range("B67").Value = 41.33761
$wsDetailed.Range("B72").Value = 36.05949
$wsDetailed.Range("B73").Value = 36.07
$wsDetailed.Range("B74").Value = 36.05962
$wsDetailed.Range("B75").Value = 36.0601
$wsDetailed.Range("B77").Value = 29.31745
$wsDetailed.Range("B78").Value = 0.51
$wsDetailed.Range("B79").Value = 6.78993
$wsDetailed.Range("B80").Value = 4.55473
$wsDetailed.Range("B81").Value = 15.96168
$wsDetailed.Range("B82").Value = 20.88137
$wsDetailed.Range("B83").Value = 1.66636
$wsDetailed.Range("B84").Value = -5.75644
$wsDetailed.Range("B85").Value = -9.914669999999999
$wsDetailed.Range("B86").Value = -6.83274
$wsDetailed.Range("B87").Value = -3.03043
$wsDetailed.Range("B88").Value = -3.09257
